# Updated Global_M2 for easier usage.
# Updates revised C2M data points for Hungary_M2 and appends two new
# monthly observations (rows 388 and 389).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the revised values in columns C:F for existing rows ---
# (All four columns C,D,E,F always carry the same value in this sheet.)
$updates = @{
    358 = 30734900000000
    359 = 31322400000000
    361 = 33495800000000
    362 = 32988400000000
    363 = 33552400000000
    364 = 33970200000000
    365 = 33881500000000
    366 = 33795100000000
    367 = 34176400000000
    368 = 34635200000000
    369 = 34976400000000
    370 = 35509700000000
    371 = 36350700000000
    372 = 37653400000000
    373 = 38869800000000
    375 = 39527800000000
    376 = 39600300000000
    377 = 39874100000000
    380 = 41209500000000
    381 = 42040500000000
    382 = 41918300000000
    387 = 40513200000000
}

foreach ($r in $updates.Keys) {
    $val = $updates[$r]
    $ws.Cells.Item($r, 3).Value = $val
    $ws.Cells.Item($r, 4).Value = $val
    $ws.Cells.Item($r, 5).Value = $val
    $ws.Cells.Item($r, 6).Value = $val
}

# --- Append two new rows (388, 389) with the same row formatting as row 387 ---
$ws.Range("A387:G387").Copy($ws.Range("A388:G389"))

$ws.Cells.Item(388, 1).Value = 44986.45833333334
$ws.Cells.Item(388, 2).Value = "ECONOMICS:HUM2"
$ws.Cells.Item(388, 3).Value = 40410200000000
$ws.Cells.Item(388, 4).Value = 40410200000000
$ws.Cells.Item(388, 5).Value = 40410200000000
$ws.Cells.Item(388, 6).Value = 40410200000000
$ws.Cells.Item(388, 7).Value = 0

$ws.Cells.Item(389, 1).Value = 45017.45833333334
$ws.Cells.Item(389, 2).Value = "ECONOMICS:HUM2"
$ws.Cells.Item(389, 3).Value = 39781400000000
$ws.Cells.Item(389, 4).Value = 39781400000000
$ws.Cells.Item(389, 5).Value = 39781400000000
$ws.Cells.Item(389, 6).Value = 39781400000000
$ws.Cells.Item(389, 7).Value = 0
